$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28, shifting existing rows 28-113 down to 29-114
$ws.Rows("28").Insert()

# Populate the newly inserted row 28 with the new record's data
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 45076
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112037
$ws.Range("G28").Value = "Cebollín"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 4000
$ws.Range("L28").Value = 4200
$ws.Range("M28").Value = 4100
$ws.Range("N28").Value = "$/paquete 36 unidades"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 114
$ws.Range("Q28").Value = 36
$ws.Range("R28").Value = "Hortaliza"
